# Update "想去人数" (F column) figures across the 展览, 演出 and 全部类型 sheets
# to reflect newly generated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    5  = 1772
    6  = 661
    7  = 318
    8  = 479
    9  = 4454
    13 = 978
    14 = 1283
    17 = 2965
    18 = 1798
    22 = 13
    25 = 298
    26 = 28
    27 = 2305
    29 = 2383
    31 = 1108
    32 = 560
    34 = 881
    35 = 411
    36 = 1095
    37 = 895
    38 = 1165
    40 = 831
    41 = 514
    42 = 356
    44 = 3470
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    5  = 23
    18 = 1
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    6  = 1772
    7  = 661
    8  = 318
    9  = 479
    10 = 4454
    15 = 1283
    16 = 2965
    18 = 1798
    25 = 13
    27 = 298
    28 = 2305
    33 = 2383
    34 = 1108
    35 = 560
    36 = 881
    37 = 1095
    38 = 895
    40 = 1165
    41 = 831
    42 = 514
    44 = 356
    48 = 3470
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
